$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format (style matching E1) to new header date cells F1:G1
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# Copy the data-row format (style matching C2, font size 12) to new data cells F2:G27
$ws.Range("C2").Copy()
$ws.Range("F2:G27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row: new measurement dates
$ws.Range("F1").Value = 45364
$ws.Range("G1").Value = 45371

# New weight measurements for weeks 2 and 3
$ws.Range("F2").Value = 20.1
$ws.Range("G2").Value = 20.5
$ws.Range("F3").Value = 21.1
$ws.Range("G3").Value = 21.3
$ws.Range("F4").Value = 22.9
$ws.Range("G4").Value = 23
$ws.Range("F5").Value = 21.8
$ws.Range("G5").Value = 22.3
$ws.Range("F6").Value = 21.6
$ws.Range("G6").Value = 23
$ws.Range("F7").Value = 19.4
$ws.Range("G7").Value = 20
$ws.Range("F8").Value = 21.2
$ws.Range("G8").Value = 21.5
$ws.Range("F9").Value = 20.9
$ws.Range("G9").Value = 21.5
$ws.Range("F10").Value = 19.5
$ws.Range("G10").Value = 20.2
$ws.Range("F11").Value = 20.4
$ws.Range("G11").Value = 21.8
$ws.Range("F12").Value = 20.4
$ws.Range("G12").Value = 22
$ws.Range("F13").Value = 21
$ws.Range("G13").Value = 22.3
$ws.Range("F14").Value = 18.7
$ws.Range("G14").Value = 20.8
$ws.Range("F15").Value = 19
$ws.Range("G15").Value = 19.5
$ws.Range("F16").Value = 23.7
$ws.Range("G16").Value = 25.2
$ws.Range("F17").Value = 20.7
$ws.Range("G17").Value = 21.9
$ws.Range("F18").Value = 20.1
$ws.Range("G18").Value = 21.3
$ws.Range("F19").Value = 20.8
$ws.Range("G19").Value = 21.6
$ws.Range("F20").Value = 19.5
$ws.Range("G20").Value = 19.8
$ws.Range("F21").Value = 20.9
$ws.Range("G21").Value = 21.5
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 19.8
$ws.Range("F23").Value = 20.1
$ws.Range("G23").Value = 21.6
$ws.Range("F24").Value = 21.4
$ws.Range("G24").Value = 23.7
$ws.Range("F25").Value = 21.7
$ws.Range("G25").Value = 22.6
$ws.Range("F26").Value = 19.9
$ws.Range("G26").Value = 19.2
$ws.Range("F27").Value = 21.1
$ws.Range("G27").Value = 22

# Update view: scroll/select to reflect the newly added columns
$ws.Range("G28").Select()
